$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the ParentID (column B) for rows 27-30: these question rows
# belong under ID 25 (row 26, "B1. Stammdaten & Betreiber"), not 26.
$ws.Range("B27:B30").Value = "25"

# Restore the view state: the sheet was scrolled down (row 7 at the top)
# with B31 as the active/selected cell.
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B31").Select()
